$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 33) continuing the existing reg_center_user rows,
# mirroring the columns/types used by the rows above it (regcntr_id, usr_id,
# lang_code, is_active, cr_by, cr_dtimes).
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Move the window/selection the way Excel would after typing the new row in,
# then clicking back up into the existing data (as captured in the saved view).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("C31").Select()
